$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge run-split paragraphs back into single runs (identical text,
#    identical formatting on both halves -- Find/Replace across the
#    split point naturally collapses them into one <w:r>).
# ---------------------------------------------------------------------

$d.Content.Find.Execute(
    "Проектът трябва да бъде приложение реализирано чрез слой за данни, слой за услуги и презентационен слой. Приложението е задължително да има Web базиран интерфейс. За приложението трябва да се използва база данни и Entity Framework. Препоръчително е използването на външни библиотеки, с помощта на които да се реализират и други допълнителни функционалности.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Проектът трябва да бъде приложение реализирано чрез слой за данни, слой за услуги и презентационен слой. Приложението е задължително да има Web базиран интерфейс. За приложението трябва да се използва база данни и Entity Framework. Препоръчително е използването на външни библиотеки, с помощта на които да се реализират и други допълнителни функционалности.",
    2) | Out-Null

$d.Content.Find.Execute(
    "По време на работата върху проекта трябва да се използва git система. Проектите се представят пред комисия и се осъществява преглед на разработката и кода /code review/ от страна на комисията.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "По време на работата върху проекта трябва да се използва git система. Проектите се представят пред комисия и се осъществява преглед на разработката и кода /code review/ от страна на комисията.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Комисията присъжда точки за проекта според неговото представяне и нивото на софтуерната разработка, като има право и да задава допълнителни контролни въпроси по своя преценка на участниците в екипа.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Комисията присъжда точки за проекта според неговото представяне и нивото на софтуерната разработка, като има право и да задава допълнителни контролни въпроси по своя преценка на участниците в екипа.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Индивидуален принос към разработката на проекта",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Индивидуален принос към разработката на проекта",
    2) | Out-Null

$d.Content.Find.Execute(
    "Предаване на проект",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Предаване на проект",
    2) | Out-Null

$d.Content.Find.Execute(
    "Копие от всички материали свързани с разработката на проекта (код, изображения, компилирано приложение, документация и др.) се качват в системата като архивиран файл от всеки един участник в проекта преди защитата. Участник без предаден проект в системата не се оценява.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Копие от всички материали свързани с разработката на проекта (код, изображения, компилирано приложение, документация и др.) се качват в системата като архивиран файл от всеки един участник в проекта преди защитата. Участник без предаден проект в системата не се оценява.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark: delete it where it currently sits
#    (after "Преподаватели от съответния център") and recreate it,
#    collapsed, right after the "Реализация на управление на
#    резервации" run inside the table.
# ---------------------------------------------------------------------

$d.Bookmarks("_GoBack").Delete()

$findRng = $d.Content.Duplicate
$findRng.Find.Execute("Реализация на управление на резервации", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$findRng.Collapse(0)
$findRng.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $findRng)
$findRng.Text = ""

# ---------------------------------------------------------------------
# 3) Highlight (yellow) the capacity-check table row: both the run
#    text and the paragraph mark need the <w:highlight w:val="yellow"/>
#    applied -- using the table Cell's Range (rather than an
#    independently constructed Range) achieves both at once.
# ---------------------------------------------------------------------

$tbl = $d.Tables(1)
$capacityCell = $tbl.Cell(7, 2)
$capacityCell.Range.Font.HighlightColorIndex = 7
